# Corrected date.getYear bad use
# Adds the missing journal entry for row 60 (08/05/2018) on the "Feuil1"
# sheet, moves the active selection down to A61, and appends the new
# shared-string description of the work done.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# New journal entry: date / activity / hours
$ws.Cells.Item(60, 1).Value = 43228
$ws.Cells.Item(60, 2).Value = "Correction du bug lié à l'ajout de transactions et amélioration de la génération de la hashmap contenant les transactions."
$ws.Cells.Item(60, 3).Value = 2

# The activity text wraps to two lines like the rows above it, so the row
# grows from the default height to 30pt.
$ws.Rows.Item(60).RowHeight = 30

# Move the active selection to A61 (the next empty row), matching the
# author's cursor position after the edit.
$ws.Range("A61").Select() | Out-Null
